# Add a third worksheet ("Ranges") that demonstrates replacing row/column
# range references (whole-row / whole-column) with explicit areas, both on
# the same sheet and cross-sheet.

$wb = $excel.ActiveWorkbook

# New sheet goes after the existing "FormulaeTypes" sheet (i.e. becomes the
# last / 3rd tab).
$formulaeTypes = $wb.Worksheets.Item($wb.Worksheets.Count)
$ranges = $wb.Worksheets.Add($null, $formulaeTypes)
$ranges.Name = "Ranges"

# Match the page setup used by the other sheets in the workbook.
$ranges.PageSetup.PaperSize = 9
$ranges.PageSetup.Orientation = 1
$ranges.PageSetup.LeftMargin = 54
$ranges.PageSetup.RightMargin = 54
$ranges.PageSetup.TopMargin = 72
$ranges.PageSetup.BottomMargin = 72
$ranges.PageSetup.HeaderMargin = 36
$ranges.PageSetup.FooterMargin = 36

# Row labels.
$ranges.Range("A2").Value = "Standard"
$ranges.Range("A3").Value = "Column"
$ranges.Range("A4").Value = "Row"

# Column headers.
$ranges.Range("B1").Value = "This sheet"
$ranges.Range("C1").Value = "Other sheet"

# Backing values summed by the "this sheet" formulas below.
$ranges.Range("F4").Value = 1
$ranges.Range("E5").Value = 1
$ranges.Range("F5").Value = 2
$ranges.Range("G5").Value = 3
$ranges.Range("F6").Value = 3

# "This sheet" column: standard area, full column, full row.
$ranges.Range("B2").Formula = "=SUM(F4:F6)"
$ranges.Range("B3").Formula = "=SUM(F:F)"
$ranges.Range("B4").Formula = "=SUM(5:5)"

# "Other sheet" column: same three flavours, referencing ValueTypes.
$ranges.Range("C2").Formula = "=SUM(ValueTypes!A3:A4)"
$ranges.Range("C3").Formula = "=SUM(ValueTypes!A:A)"
$ranges.Range("C4").Formula = "=SUM(ValueTypes!4:4)"

# Make the new sheet the active / selected tab, with C5 selected.
[void]$ranges.Range("C5").Select()
